$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2:B9 label text (Variação 2021/2012 -> Variação 2022/2013)
$ws.Range("B2:B9").Value = "Variação 2022/2013"

# Row 2: Mato Grosso
$ws.Range("C2").Value = 31.57612616475882

# Row 3: Tocantins
$ws.Range("C3").Value = 30.25169296719681

# Row 4: was Santa Catarina -> now Maranhão
$ws.Range("A4").Value = "Maranhão"
$ws.Range("C4").Value = 26.68172809047942

# Row 5: was Roraima -> now Mato Grosso do Sul
$ws.Range("A5").Value = "Mato Grosso do Sul"
$ws.Range("C5").Value = 22.47693215339233

# Row 6: was Maranhão -> now Santa Catarina
$ws.Range("A6").Value = "Santa Catarina"
$ws.Range("C6").Value = 20.57128073428023

# Row 7: Acre
$ws.Range("C7").Value = 20.37343691964839

# Row 8: Sergipe
$ws.Range("C8").Value = 2.799581048610689
$ws.Range("D8").Value = "23º"

# Row 9: Nordeste
$ws.Range("C9").Value = 9.525323439430643

# Row 10: Brasil - remove entire row
$ws.Rows.Item(10).Delete()
